# Estadisticos Segundo Parcial 26 Mayo
#
# 1. Calificaciones sheet: fill in the "2P" (segundo parcial) grade (column H)
#    for each student in the first subject block, and refresh the "Final"
#    grade (column T) that depends on it.
# 2. Asistencias sheet: refresh the attendance percentages for 2P (column H)
#    and 3P (column N) which mirror each other for the first subject block.
# 3. Totales sheet: refresh the aggregate stats (Aprobados/Reprobados/
#    Por_Apro/Por_Repro/Promedio) for the first subject/teacher row.
# 4. Rescatables sheet: drop the second "rescatable" student row (PUGA ROMERO
#    HILEN ALELI) - only CHONKOA SANDOVAL ABRIL remains.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Calificaciones — column H (2P) and column T (Final) for rows 4-28
# ---------------------------------------------------------------------
$wsCal = $wb.Worksheets.Item("Calificaciones")

$calValues = @{
  "H4" = 7;  "H5" = 9;  "T5" = 9;  "H6" = 9;  "H7" = 5;
  "H8" = 9;  "T8" = 9;  "H9" = 9;  "T9" = 9;  "H10" = 9;
  "H11" = 8; "T11" = 7; "H12" = 9; "T12" = 9; "H13" = 8;
  "T13" = 8; "H14" = 9; "T14" = 9; "H15" = 9; "H16" = 9;
  "T16" = 9; "H17" = 9; "T17" = 9; "H18" = 8; "T18" = 8;
  "H19" = 9; "T19" = 9; "H20" = 8; "T20" = 7; "H21" = 8;
  "T21" = 8; "H22" = 9; "T22" = 9; "H23" = 9; "T23" = 9;
  "H24" = 7; "H25" = 10; "T25" = 10; "H26" = 9; "T26" = 9;
  "H27" = 7; "H28" = 8; "T28" = 8
}

foreach ($ref in $calValues.Keys) {
  $wsCal.Range($ref).Value = $calValues[$ref]
}

# ---------------------------------------------------------------------
# 2. Asistencias — columns H (2P) and N (3P) mirror each other, rows 4-28
# ---------------------------------------------------------------------
$wsAsis = $wb.Worksheets.Item("Asistencias")

$asisValues = @{
  4  = 97.5;
  5  = 96.2;
  7  = 91.09999999999999;
  8  = 97.5;
  9  = 97.5;
  10 = 98.7;
  11 = 87.3;
  12 = 98.7;
  13 = 94.90000000000001;
  14 = 97.5;
  15 = 97.5;
  16 = 98.7;
  17 = 87.3;
  18 = 97.5;
  19 = 97.5;
  20 = 82.3;
  21 = 98.7;
  22 = 97.5;
  23 = 97.5;
  24 = 98.7;
  26 = 98.7;
  27 = 93.7;
  28 = 87.3
}

foreach ($row in $asisValues.Keys) {
  $val = $asisValues[$row]
  $wsAsis.Range("H$row").Value = $val
  $wsAsis.Range("N$row").Value = $val
}

# ---------------------------------------------------------------------
# 3. Totales — first subject/teacher row (row 2) aggregate recompute
# ---------------------------------------------------------------------
$wsTot = $wb.Worksheets.Item("Totales")

$wsTot.Range("D2").Value = 24
$wsTot.Range("E2").Value = 1
$wsTot.Range("F2").Value = 96
$wsTot.Range("G2").Value = 4
$wsTot.Range("H2").Value = 8.4

# ---------------------------------------------------------------------
# 4. Rescatables — remove the PUGA ROMERO HILEN ALELI row (row 3)
# ---------------------------------------------------------------------
$wsResc = $wb.Worksheets.Item("Rescatables")
$wsResc.Rows.Item(3).Delete()
